$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.934604333333334
$ws.Range("H2").Value = 17.803813
$ws.Range("I2").Value = 0.3081877218757661
$ws.Range("J2").Value = 0.3081877218757661
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.51723133333333
$ws.Range("N2").Value = 34.551694
$ws.Range("Q2").Value = 68.35021097880245
$ws.Range("R2").Value = 615.151898809222
$ws.Range("S2").Value = 0.3081877218757661
$ws.Range("T2").Value = 0.3081877218757661

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.720664
$ws.Range("H3").Value = 23.161992
$ws.Range("I3").Value = 0.4009389195777736
$ws.Range("J3").Value = 0.4009389195777736
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 11.51723133333333
$ws.Range("N3").Value = 34.551694
$ws.Range("Q3").Value = 88.92067333493867
$ws.Range("R3").Value = 800.286060014448
$ws.Range("S3").Value = 0.4009389195777736
$ws.Range("T3").Value = 0.4009389195777736

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.601191
$ws.Range("H4").Value = 16.803573
$ws.Range("I4").Value = 0.2908733585464604
$ws.Range("J4").Value = 0.2908733585464603
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 11.51723133333333
$ws.Range("N4").Value = 34.551694
$ws.Range("Q4").Value = 64.51021248918467
$ws.Range("R4").Value = 580.5919124026619
$ws.Range("S4").Value = 0.2908733585464604
$ws.Range("T4").Value = 0.2908733585464603
